$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell values (columns A, B, D, E for rows 1-4)
$ws.Range("B1").Value = 86

$ws.Range("A2").Value = 43
$ws.Range("B2").Value = 21
$ws.Range("D2").Value = 99
$ws.Range("E2").Value = 17

$ws.Range("B3").Value = 35

$ws.Range("A4").Value = 65
$ws.Range("B4").Value = 4
$ws.Range("D4").Value = 345
$ws.Range("E4").Value = 80

# Add new row 5
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 8
$ws.Range("D5").Value = 22
$ws.Range("E5").Value = 19

# Update the sheet view selection from E9 to F9
$ws.Range("F9").Select()
